$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 1, 0.98),
    @(0.39, 0.311, 0.35),
    @(0.425, 0.423, 0.28),
    @(0.003, 0.117, 0.25),
    @(0.013, 0, 0.27),
    @(0.003, 0, 0.22),
    @(0.228, 0.228, 0.19),
    @(0.002, 0, 0.2),
    @(0.347, 0.326, 0.2),
    @(0.182, 0.171, 0.18),
    @(1, 1, 0.94),
    @(0.002, 0.039, 0.2),
    @(0.049, 0.172, 0.2),
    @(0, 0, 0.14),
    @(0.01, 0.181, 0.15),
    @(0.103, 0.042, 0.15),
    @(0.012, 0.009, 0.11),
    @(0.03, 0.105, 0.12),
    @(0.017, 0, 0.77)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
